# [Kadastro App] Yeni kayit eklendi: 1
# Appends a new record both to the master "Kayitlar" log sheet and to the
# district-specific sheet ("Anamur") matching the record's Birim value.
# Every column in these tables is plain text (even the numeric-looking
# ones like "Kayit No" / "Parsel Sayisi"), so each value is written through
# a helper that forces the Text number format before assigning it, then
# restores the default ("Normal") style so no stray formatting is left on
# the cell.

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

function Add-KayitRow($ws, $row, $kayitNo, $tarih, $birim, $parsel, $is, $personel) {
    Set-TextCell $ws.Cells.Item($row, 1) $kayitNo
    Set-TextCell $ws.Cells.Item($row, 2) $tarih
    Set-TextCell $ws.Cells.Item($row, 3) $birim
    Set-TextCell $ws.Cells.Item($row, 4) $parsel
    Set-TextCell $ws.Cells.Item($row, 5) $is
    Set-TextCell $ws.Cells.Item($row, 6) $personel
}

$wb = $excel.ActiveWorkbook

$kayitNo  = "1"
$tarih    = "2025-09-08"
$birim    = "Anamur"
$parsel   = "50"
$is       = "18-UYG."
$personel = "EMİNE ALANLI KIRCILI (K.Mühendisi), HALİL ÇETİNKAYA (K.Teknisyeni)"

# --- Kayitlar (master list) sheet: append the new record as the next row ---
$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
$lastRowKayitlar = $wsKayitlar.Cells.Item($wsKayitlar.Rows.Count, 1).End(-4162).Row
$newRowKayitlar = $lastRowKayitlar + 1
Add-KayitRow $wsKayitlar $newRowKayitlar $kayitNo $tarih $birim $parsel $is $personel

# --- District sheet (named after Birim, here "Anamur"): same new record ---
$wsBirim = $wb.Worksheets.Item($birim)
$lastRowBirim = $wsBirim.Cells.Item($wsBirim.Rows.Count, 1).End(-4162).Row
$newRowBirim = $lastRowBirim + 1
Add-KayitRow $wsBirim $newRowBirim $kayitNo $tarih $birim $parsel $is $personel
